$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Control 0
$ws.Range("D2").Value = [double]"7.207512075506461E-10"
$ws.Range("E2").Value = [double]"7.207512075506461E-10"

# Row 3 - Control 6
$ws.Range("D3").Value = [double]"3.341207309038588E-05"
$ws.Range("E3").Value = [double]"3.341207309038588E-05"

# Row 4 - Control 9
$ws.Range("C4").Value = $true
$ws.Range("D4").Value = [double]"0.1225870191204658"
$ws.Range("E4").Value = [double]"0.1225870191204658"

# Row 5 - Control 24
$ws.Range("D5").Value = [double]"0.9999873609829134"
$ws.Range("E5").Value = [double]"0.9999873609829134"

# Row 6 - Control 32
$ws.Range("D6").Value = [double]"0.9999999998212981"
$ws.Range("E6").Value = [double]"0.9999999998212981"

# Row 8 - MDD 12
$ws.Range("D8").Value = [double]"0.5278789314385657"
$ws.Range("E8").Value = [double]"0.4721210685614343"

# Row 9 - MDD 53
$ws.Range("D9").Value = [double]"0.999999970955868"
$ws.Range("E9").Value = [double]"2.904413198834277E-08"

# Row 10 - MDD 29
$ws.Range("D10").Value = [double]"5.428819304239419E-22"

# Row 11 - MDD 55
$ws.Range("D11").Value = [double]"1.745251430907637E-25"
$ws.Range("F11").Value = [double]"14.04665946960449"
$ws.Range("G11").Value = [double]"0.6"

Write-Output "done"
